# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.889.76"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "3.406.11"
$ws.Range("E3").Value = "  -3.09%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "405.61"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "134.57"
$ws.Range("E6").Value = "  +3.36%  "

$ws.Range("E7").Value = "  -1.63%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.94%  "

$ws.Range("E10").Value = "  -6.79%  "

$ws.Range("D11").Value = "42.57"
$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("E13").Value = "  -4.21%  "

$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").Value = "3.395.69"
$ws.Range("E15").Value = "  -2.39%  "

$ws.Range("D16").Value = "61.892.40"
$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("E17").Value = "  -3.59%  "

$ws.Range("D18").Value = "11.01"
$ws.Range("E18").Value = "  -3.13%  "

$ws.Range("D19").Value = "0.0000129"
$ws.Range("E19").Value = "  -6.32%  "

$ws.Range("E20").Value = "  -5.54%  "

$ws.Range("D21").Value = "84.39"
$ws.Range("E21").Value = "  +2.53%  "

$ws.Range("D22").Value = "311.76"
$ws.Range("E22").Value = "  -0.90%  "

$ws.Range("D23").Value = "12.89"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("E25").Value = "  +8.88%  "

$ws.Range("D26").Value = "29.52"
$ws.Range("E26").Value = "  -3.71%  "

$ws.Range("D27").Value = "8.14"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").Value = "2.81"
$ws.Range("E28").Value = "  +4.53%  "

$ws.Range("D29").Value = "7.60"
$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("E30").Value = "  -4.66%  "

$ws.Range("E31").Value = "  -2.94%  "

$ws.Range("D32").Value = "42.43"
$ws.Range("E32").Value = "  -2.95%  "

$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").Value = "11.30"

$ws.Range("D35").Value = "0.0481"
$ws.Range("E35").Value = "  -3.00%  "

$ws.Range("D36").Value = "51.64"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  -5.28%  "

$ws.Range("E39").Value = "  -3.62%  "

$ws.Range("D40").Value = "0.303"
$ws.Range("E40").Value = "  +4.94%  "

$ws.Range("D41").Value = "137.13"
$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").Value = "1.97"
$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("D44").Value = "4.01"
$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("D45").Value = "16.65"
$ws.Range("E45").Value = "  -6.53%  "

$ws.Range("D47").Value = "21.25"
$ws.Range("E47").Value = "  -5.67%  "

$ws.Range("D48").Value = "2.115.03"
$ws.Range("E48").Value = "  -5.07%  "

$ws.Range("E49").Value = "  -2.72%  "

$ws.Range("D50").Value = "1.90"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("E51").Value = "  +15.05%  "

